# edges.xlsx update: add more edges (rows) to the "edges" sheet, per the
# commit "Actualización edges, añadiendo aristas" — adds rows 8-29 and
# corrects the B/C values on existing rows 6-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (row, A, B, C) for every row that needs to exist/ change, rows 6-29.
$rows = @(
    @(6, 2, 9, 82.18),
    @(7, 2, 10, 106.83),
    @(8, 2, 3, 88.83),
    @(9, 2, 39, 100.35),
    @(10, 2, 40, 191.89),
    @(11, 2, 34, 247.19),
    @(12, 2, 12, 176.12),
    @(13, 2, 14, 208.72),
    @(14, 2, 5, 145.22),
    @(15, 3, 9, 147.52),
    @(16, 3, 10, 122.61),
    @(17, 3, 11, 134.4),
    @(18, 3, 31, 323.96),
    @(19, 3, 19, 192.28),
    @(20, 3, 12, 130.83),
    @(21, 3, 14, 142.3),
    @(22, 3, 8, 145.69),
    @(23, 3, 4, 62.43),
    @(24, 3, 40, 141.4),
    @(25, 3, 39, 135.25),
    @(26, 3, 41, 175.3),
    @(27, 3, 43, 248.48),
    @(28, 3, 5, 67.61),
    @(29, 3, 18, 178.63)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Page setup - force orientation=portrait so <pageSetup> gets serialized.
$ws.PageSetup.Orientation = 1

# Move the selection to the new first empty row below the data (C30),
# mirroring the author's cursor position after entering the last edge.
$ws.Range("C30").Select() | Out-Null
